# Auto update Excel log
# Appends new sensor-log rows to the "Proximity" and "Camera" sheets,
# mirroring a fresh batch of door-sensor events captured at 15:16.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Proximity sheet: three new rows (50-52) for the Living Room Main Door
# ---------------------------------------------------------------------
$proximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "15:16:12", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:16:14", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:16:15", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
)

$startRow = 50
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $row = $startRow + $i
    $values = $proximityRows[$i]

    # Dates like "2026-02-01" would otherwise be auto-recognised as a date
    # serial by Excel's input parser. Prefix with an apostrophe to force
    # plain text entry (exactly like typing '2026-02-01 into the UI), then
    # clear the resulting "quote prefix" cell style so formatting stays
    # identical to a cell that was never touched.
    $dateCell = $proximity.Cells.Item($row, 1)
    $dateCell.Value = "'" + $values[0]
    $dateCell.Style = "Normal"

    $proximity.Cells.Item($row, 2).Value = $values[1]
    $proximity.Cells.Item($row, 3).Value = $values[2]
    $proximity.Cells.Item($row, 4).Value = $values[3]
    $proximity.Cells.Item($row, 5).Value = $values[4]
    $proximity.Cells.Item($row, 6).Value = $values[5]
}

# ---------------------------------------------------------------------
# Camera sheet: one new row (17) - image capture tied to the same event
# ---------------------------------------------------------------------
$camera = $wb.Worksheets.Item("Camera")

$cameraRow = 17
$cameraValues = @("2026-02-01", "15:16:15", "15:00", "Living Room Main Door", "Image Captured", "Active")

$cameraDateCell = $camera.Cells.Item($cameraRow, 1)
$cameraDateCell.Value = "'" + $cameraValues[0]
$cameraDateCell.Style = "Normal"

$camera.Cells.Item($cameraRow, 2).Value = $cameraValues[1]
$camera.Cells.Item($cameraRow, 3).Value = $cameraValues[2]
$camera.Cells.Item($cameraRow, 4).Value = $cameraValues[3]
$camera.Cells.Item($cameraRow, 5).Value = $cameraValues[4]
$camera.Cells.Item($cameraRow, 6).Value = $cameraValues[5]
